$d = $word.ActiveDocument

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$oldColor = RGBVal 0x11 0x11 0x11
$newColor = RGBVal 0x00 0x00 0x00

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Text = ""
$find.Font.Color = $oldColor
$find.Replacement.Text = ""
$find.Replacement.Font.Color = $newColor

$find.Execute(
    [ref]"",        # FindText
    [ref]$false,    # MatchCase
    [ref]$false,    # MatchWholeWord
    [ref]$false,    # MatchWildcards
    [ref]$false,    # MatchSoundsLike
    [ref]$false,    # MatchAllWordForms
    [ref]$true,     # Forward
    [ref]1,         # Wrap (wdFindContinue)
    [ref]$true,     # Format
    [ref]"",        # ReplaceWith
    [ref]2          # Replace (wdReplaceAll)
)

# Find/Replace skips runs with empty text (e.g. the blank first header
# cell), so sweep every table cell directly and fix any leftover runs
# that still carry the old color.
foreach ($table in $d.Tables) {
    foreach ($cell in $table.Range.Cells) {
        $cellRange = $cell.Range
        if ($cellRange.Font.Color -eq $oldColor) {
            $cellRange.Font.Color = $newColor
        }
    }
}

